$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value2 = "Dataset"
$ws.Range("B1").Value2 = "Number of Examples (In-Depth Runs)"
$ws.Range("C1").Value2 = "Number of Examples (Optimized Single Run)"
$ws.Range("D1").Value2 = "Average Training Set Accuracy (%)"
$ws.Range("E1").Value2 = "Average Testing Set Accuracy (%)"
$ws.Range("F1").Value2 = "Average Time Elapsed (s)"
$ws.Range("G1").Value2 = "Max Depth"
$ws.Range("H1").Value2 = "Min Examples"
$ws.Range("I1").Value2 = "Min Proportion"
$ws.Range("J1").Value2 = "Num Trees"

# ---- Data rows (already in final, alphabetically-sorted-by-Dataset order) ----
# Row 2: Adult
$ws.Range("A2").Value2 = "Adult"
$ws.Range("B2").Value2 = 500
$ws.Range("C2").Value2 = 2500
$ws.Range("D2").Value2 = 83.51
$ws.Range("E2").Value2 = 81.400000000000006
$ws.Range("F2").Value2 = 100.68
$ws.Range("G2").Value2 = 3
$ws.Range("H2").Value2 = "N/A"
$ws.Range("I2").Value2 = "N/A"
$ws.Range("J2").Value2 = 100

# Row 3: Blobs
$ws.Range("A3").Value2 = "Blobs"
$ws.Range("B3").Value2 = 600
$ws.Range("C3").Value2 = 600
$ws.Range("D3").Value2 = 95.88
$ws.Range("E3").Value2 = 94
$ws.Range("F3").Value2 = 3.61
$ws.Range("G3").Value2 = 4
$ws.Range("H3").Value2 = 80
$ws.Range("I3").Value2 = 0.75
$ws.Range("J3").Value2 = 100

# Row 4: Digits
$ws.Range("A4").Value2 = "Digits"
$ws.Range("B4").Value2 = 250
$ws.Range("C4").Value2 = 250
$ws.Range("D4").Value2 = 100
$ws.Range("E4").Value2 = 65.599999999999994
$ws.Range("F4").Value2 = 35.630000000000003
$ws.Range("G4").Value2 = "N/A"
$ws.Range("H4").Value2 = "N/A"
$ws.Range("I4").Value2 = "N/A"
$ws.Range("J4").Value2 = 150

# Row 5: Letters
$ws.Range("A5").Value2 = "Letters"
$ws.Range("B5").Value2 = 500
$ws.Range("C5").Value2 = 500
$ws.Range("D5").Value2 = 99.9
$ws.Range("E5").Value2 = 67.8
$ws.Range("F5").Value2 = 47.29
$ws.Range("G5").Value2 = 9
$ws.Range("H5").Value2 = "N/A"
$ws.Range("I5").Value2 = "N/A"
$ws.Range("J5").Value2 = 100

# Row 6: Spirals
$ws.Range("A6").Value2 = "Spirals"
$ws.Range("B6").Value2 = 0
$ws.Range("C6").Value2 = 1000
$ws.Range("D6").Value2 = 100
$ws.Range("E6").Value2 = 96.7
$ws.Range("F6").Value2 = 46.39
$ws.Range("G6").Value2 = "N/A"
$ws.Range("H6").Value2 = "N/A"
$ws.Range("I6").Value2 = "N/A"
$ws.Range("J6").Value2 = 100

# Row 7: Zoo
$ws.Range("A7").Value2 = "Zoo"
$ws.Range("B7").Value2 = 0
$ws.Range("C7").Value2 = 101
$ws.Range("D7").Value2 = 100
$ws.Range("E7").Value2 = 97
$ws.Range("F7").Value2 = 0.18
$ws.Range("G7").Value2 = "N/A"
$ws.Range("H7").Value2 = "N/A"
$ws.Range("I7").Value2 = "N/A"
$ws.Range("J7").Value2 = 150

# ---- Number formatting (creates the shared "0.00" style used by all data cells) ----
$ws.Range("B2:J7").NumberFormat = "0.00"
$ws.Range("J8").NumberFormat = "0.00"

# ---- Column widths (best-fit approximation) ----
$ws.Columns("A").ColumnWidth = 7.33203125
$ws.Columns("B").ColumnWidth = 31.44140625
$ws.Columns("C").ColumnWidth = 37
$ws.Columns("D").ColumnWidth = 28.6640625
$ws.Columns("E").ColumnWidth = 28
$ws.Columns("F").ColumnWidth = 21.33203125
$ws.Columns("G").ColumnWidth = 9.88671875
$ws.Columns("H").ColumnWidth = 12.21875
$ws.Columns("I").ColumnWidth = 13.44140625
$ws.Columns("J").ColumnWidth = 9.77734375

# ---- Sort bookkeeping (matches the sortState Excel records after sorting by Dataset) ----
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A7"))
$sortObj.SetRange($ws.Range("A1:J7"))
$sortObj.Header = 1
$sortObj.Apply()

# ---- Page setup ----
$ws.PageSetup.Orientation = 1

# ---- Selection ----
$ws.Range("G7").Select()
